$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph about "Die Schopfung" right before the existing
#    "Une folie" paragraph (5 April 1802 entry), splitting what used to be a
#    single paragraph into two - the first one keeps the "5 April 1802"
#    heading and gains the new Haydn sentence, the second one keeps the
#    original "Une folie" sentence untouched.
# ---------------------------------------------------------------------------

$titleText = "Die Sch" + [char]0x00F6 + "pfung"
$bodyText  = " by Franz Joseph Haydn (70) is performed in Engelhardt House, St. Petersburg to inaugurate the St. Petersburg Philharmonic Society."

$target = $d.Content
$found = $target.Find.Execute("Une folie", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Une folie' anchor text"
}

$splitPoint = $target.Start

# Split the paragraph: this creates a new empty paragraph boundary right
# before "Une folie", re-using the paragraph's existing pPr (color only, no
# bold heading run) - exactly matching the target structure.
$target.InsertParagraphBefore()

# Insert the italic title run at the end of the now-split first paragraph
# (right before the new paragraph mark).
$titleIns = $d.Range($splitPoint, $splitPoint)
$titleIns.InsertAfter($titleText)
$titleEnd = $splitPoint + $titleText.Length
$titleRng = $d.Range($splitPoint, $titleEnd)
$titleRng.Italic = 1
$titleRng.Font.Color = 0

# Insert the regular (non-italic) body run right after the title run.
$bodyIns = $d.Range($titleEnd, $titleEnd)
$bodyIns.InsertAfter($bodyText)
$bodyEnd = $titleEnd + $bodyText.Length
$bodyRng = $d.Range($titleEnd, $bodyEnd)
$bodyRng.Font.Color = 0

# ---------------------------------------------------------------------------
# 2) Update the copyright line: "2004-2013" -> "2004-2016".
#    Runs "c" / "2004-2013" / " " / "Paul Scharfenberger " have no explicit
#    formatting, so a plain text edit would cause the whole paragraph to be
#    re-coalesced into a single run. Wrapping the untouched neighbour runs
#    in a temporary bookmark forces the engine to keep run boundaries intact;
#    the bookmark itself is removed again once the edit is done.
# ---------------------------------------------------------------------------

$yearRng = $d.Content
$foundYear = $yearRng.Find.Execute("2004-2013", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundYear) {
    throw "Could not find '2004-2013' text"
}

$afterYearRng = $d.Range($yearRng.End, $yearRng.End)
$foundAfter = $afterYearRng.Find.Execute("Paul Scharfenberger ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundAfter) {
    throw "Could not find 'Paul Scharfenberger ' text"
}

$bmAfter = $d.Bookmarks.Add("zzEditWallAfter", $afterYearRng)

$yearRng.Text = "2004-2016"

$bmAfterObj = $d.Bookmarks("zzEditWallAfter")
$bmAfterObj.Delete()

# ---------------------------------------------------------------------------
# 3) Update the date line: "6 September 2013" -> "6 July 2016"
#    ("September" + " 2013" runs collapse into a single "July 2016" run,
#    while the leading "6 " run is left untouched/separate).
# ---------------------------------------------------------------------------

$dateRng = $d.Content
$foundDate = $dateRng.Find.Execute("September 2013", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundDate) {
    throw "Could not find 'September 2013' text"
}

$sixRng = $d.Range($dateRng.Start - 2, $dateRng.Start)

$bmSix = $d.Bookmarks.Add("zzEditWallSix", $sixRng)

$dateRng.Text = "July 2016"

$bmSixObj = $d.Bookmarks("zzEditWallSix")
$bmSixObj.Delete()
